{"js": "// Find the paragraph that currently reads \"Biblioteca AWT e Swing de java\"\n// and replace it with the expanded set of paragraphs described by the diff:\n//   1. \"Biblioteca AWT e Swing de java\" (same text, but \"java\" wrapped in\n//      proofErr spell-check tags, split across two runs)\n//   2. an empty (blank) centered paragraph\n//   3. \"heran\u00e7a\"\n//   4. \"Sobrecarga vs Sobrescrita\"\n//   5. \"Class super\"\n//   6. an empty (blank) centered paragraph\n//   7. \"Jva bins\"\n//   8. \"Atributos privados\"\n//   9. \"Polimorfismo = mesmo c\u00f3digo ter v\u00e1rios comportamentos.\" (this last\n//      paragraph keeps the trailing _GoBack bookmark that was on the\n//      original paragraph, and is NOT centered, matching the diff).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Biblioteca AWT e Swing de java\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Could not find paragraph 'Biblioteca AWT e Swing de java'\");\n}\n\nconst newBodyXml =\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Biblioteca AWT e Swing de </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>java</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t>heran\\u00e7a</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t xml:space=\"preserve\">Sobrecarga </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:proofErr w:type=\"gramStart\"/>' +\n  '<w:r><w:t>vs</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> Sobrescrita</w:t></w:r>' +\n  '<w:proofErr w:type=\"gramEnd\"/>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>Class</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>super</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>Jva</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>bins</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '</w:p>' +\n  '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n  '<w:r><w:t>Atributos privados</w:t></w:r>' +\n  '</w:p>' +\n  '<w:p>' +\n  '<w:r><w:t>P</w:t></w:r>' +\n  '<w:r><w:t>olimorfismo</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\"> = mesmo c\\u00f3digo ter v\\u00e1rios comportamentos.</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>';\n\nconst flatOpc =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<?mso-application progid=\"Word.Document\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newBodyXml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ntarget.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Locate the paragraph that currently reads \"Biblioteca AWT e Swing de java\"\n# and replace its whole range (text + paragraph mark) with the expanded set\n# of paragraphs described by the diff:\n#   1. \"Biblioteca AWT e Swing de java\" (same text, \"java\" now wrapped in\n#      proofErr spell-check tags, split across two runs)\n#   2. an empty (blank) centered paragraph\n#   3. \"heran\u00e7a\"\n#   4. \"Sobrecarga vs Sobrescrita\"\n#   5. \"Class super\"\n#   6. an empty (blank) centered paragraph\n#   7. \"Jva bins\"\n#   8. \"Atributos privados\"\n#   9. \"Polimorfismo = mesmo c\u00f3digo ter v\u00e1rios comportamentos.\" (this last\n#      paragraph keeps the trailing _GoBack bookmark that used to sit on the\n#      original paragraph, and is NOT centered, matching the diff).\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Biblioteca AWT e Swing de java*\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Could not find paragraph 'Biblioteca AWT e Swing de java'\"\n}\n\n# Use the paragraph's own Range (covers the text AND its paragraph mark) so\n# that InsertXML performs a full structural replace of the paragraph rather\n# than just swapping the text within it.\n$rng = $target.Range\n\n$newBodyXml = '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Biblioteca AWT e Swing de </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>java</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r><w:t>heran\u00e7a</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r><w:t xml:space=\"preserve\">Sobrecarga </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>vs</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Sobrescrita</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Class</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>super</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr></w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Jva</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>bins</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>' +\n    '<w:p><w:pPr><w:jc w:val=\"center\"/></w:pPr>' +\n    '<w:r><w:t>Atributos privados</w:t></w:r>' +\n    '</w:p>' +\n    '<w:p>' +\n    '<w:r><w:t>P</w:t></w:r>' +\n    '<w:r><w:t>olimorfismo</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> = mesmo c\u00f3digo ter v\u00e1rios comportamentos.</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '</w:p>'\n\n$flatOpc = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n    '<?mso-application progid=\"Word.Document\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $newBodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$rng.InsertXML($flatOpc)\n"}
